$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Headers (write B3 first so "Sensor output" lands at shared-string index 0)
$ws.Range("B3").Value = "Sensor output"
$ws.Range("A3").Value = "Object Distance (cm)"

# Calibration data (rows 5-11)
$ws.Range("A5").Value = 20
$ws.Range("B5").Value = 494
$ws.Range("A6").Value = 30
$ws.Range("B6").Value = 375
$ws.Range("A7").Value = 40
$ws.Range("B7").Value = 294
$ws.Range("A8").Value = 50
$ws.Range("B8").Value = 236
$ws.Range("A9").Value = 60
$ws.Range("B9").Value = 200
$ws.Range("A10").Value = 70
$ws.Range("B10").Value = 171
$ws.Range("A11").Value = 80
$ws.Range("B11").Value = 150

# Test section
$ws.Range("A18").Value = "Test of points not used for calibration"
$ws.Range("A19").Value = 25
$ws.Range("B19").Value = 416
$ws.Range("A20").Value = 35
$ws.Range("B20").Value = 332
$ws.Range("A21").Value = 45
$ws.Range("B21").Value = 250
$ws.Range("A22").Value = 55
$ws.Range("B22").Value = 226
$ws.Range("A23").Value = 65
$ws.Range("B23").Value = 194
$ws.Range("A24").Value = 75
$ws.Range("B24").Value = 160

# Column width for column A
$ws.Columns.Item(1).ColumnWidth = 14.42578125

# View settings
$ws.Range("A18").Select()

$chartObj = $ws.Shapes.AddChart2(-1, 4)
$chartObj.Chart.SetSourceData($ws.Range("A4:B11"))
$chartObj.Chart.ChartType = 4
